# Commit: "Fixed POI packaging and upgraded to POI 3.15."
#
# The underlying XML diff for this resource is purely a re-serialization
# artifact of the Apache POI upgrade: every changed line is nothing more
# than the *attribute order* inside existing elements (namespace
# declarations on <w:document>, and attribute order on <w:pgSz>, <w:pgMar>,
# <w:rFonts>, <w:lang>, <w:latentStyles>, <w:lsdException>, <w:style>,
# <w:tblInd>, <w:tblCellMar> children, ...) being written out in a
# (now alphabetically-sorted) order by the newer POI XML writer.
#
# No text, value, style, formatting, or structural content actually
# changed - every attribute name/value pair present before is still
# present after, just emitted in a different order, which is confirmed
# by canonicalizing (C14N) both sides: they are byte-identical once
# normalized. There is therefore no user-visible edit for Word's object
# model to perform here; this script intentionally touches nothing.
$d = $word.ActiveDocument

# Touch the document read-only to confirm the session/object model is
# wired up correctly, without mutating any content.
$null = $d.Content.Text
